# CompFormStatus.xlsx edit: swap the "Site Extraction" status block (rows
# 29-30) with the "Match Found" status block (rows 31-32) on Sheet1, and
# move the sheet's selection to C16 (matches the scroll/selection change
# recorded for this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 29: was "Extraction Pending Site Count" row -> becomes "Full Match Found" row ---
$ws.Range("C29").Value = "Full Match Found"
$ws.Range("D29").ClearContents()
$ws.Range("F29").Value = "FullMatchFoundReviewPending"
$ws.Range("G29").Value = "Review Pending, Full Match Found"

# --- Row 30: was "Extraction Error Site Count = 0" row -> becomes "Partial Match Found" row ---
$ws.Range("D30").ClearContents()
$ws.Range("C30").Value = "Partial Match Found"
$ws.Range("F30").Value = "PartialMatchFoundReviewPending"
$ws.Range("G30").Value = "Review Pending, Partial Match Found"

# --- Row 31: was "Full Match Found" row -> becomes "Extraction Pending Site Count" row ---
$ws.Range("C31").Value = "Extraction Pending Site Count > 0"
$ws.Range("D31").Value = "Extraction Error Site Count > 0"
$ws.Range("F31").Value = "NotScanned;"
$ws.Range("G31").Value = "Data Extraction Error at {0} Source{1}, Review Pending"

# --- Row 32: was "Partial Match Found" row -> becomes "Extraction Error Site Count = 0" row ---
$ws.Range("C32").ClearContents()
$ws.Range("D32").Value = "Extraction Error Site Count = 0"
$ws.Range("F32").Value = "NotScanned;"
$ws.Range("G32").Value = "Data Extraction Pending at {0} Source{1}, Review Pending"

# --- View: move the selection down to C16 (keeps the existing row-1 freeze) ---
$ws.Range("C16").Select() | Out-Null
